$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, $val) {
    # Force the value to be stored as literal text, even when it
    # looks numeric (e.g. "590.46" or "0.0000260"), matching the
    # inline-string cells in the source data. Reset to the Normal
    # style afterwards so no stray number-format is left behind.
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "67.517.39"
$ws.Range("E2").Value = "  -0.43%  "
$ws.Range("D3").Value = "3.721.55"
$ws.Range("E3").Value = "  -2.13%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("E5").Value = "  -1.44%  "
$ws.Range("E6").Value = "  -2.21%  "
$ws.Range("D7").Value = "3.719.22"
$ws.Range("E7").Value = "  -2.25%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -2.29%  "
$ws.Range("E10").Value = "  -3.96%  "
$ws.Range("E12").Value = "  -2.86%  "
$ws.Range("E13").Value = "  -5.40%  "
$ws.Range("E14").Value = "  -2.49%  "
$ws.Range("D15").Value = "4.349.53"
$ws.Range("E15").Value = "  -2.01%  "
$ws.Range("D16").Value = "3.720.81"
$ws.Range("E16").Value = "  -2.40%  "
$ws.Range("D17").Value = "67.514.91"
$ws.Range("E17").Value = "  -0.52%  "
$ws.Range("E18").Value = "  -0.83%  "
$ws.Range("E19").Value = "  -5.84%  "
$ws.Range("E20").Value = "  -0.46%  "
$ws.Range("E21").Value = "  -1.67%  "
$ws.Range("E22").Value = "  -0.69%  "
$ws.Range("E23").Value = "  -4.88%  "
$ws.Range("E24").Value = "  -0.77%  "
$ws.Range("E25").Value = "  -11.35%  "
$ws.Range("E26").Value = "  -6.81%  "
$ws.Range("E27").Value = "  -2.02%  "
$ws.Range("E28").Value = "  -2.70%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").Value = "3.868.23"
$ws.Range("E30").Value = "  -1.90%  "
$ws.Range("E31").Value = "  -5.70%  "
$ws.Range("E32").Value = "  -5.51%  "
$ws.Range("E33").Value = "  -2.94%  "
$ws.Range("E34").Value = "  -4.22%  "
$ws.Range("E35").Value = "  -3.22%  "
$ws.Range("D36").Value = "3.674.82"
$ws.Range("E36").Value = "  -2.36%  "
$ws.Range("E37").Value = "  -5.62%  "
$ws.Range("E38").Value = "  -9.22%  "
$ws.Range("E39").Value = "  -1.61%  "
$ws.Range("E40").Value = "  -2.18%  "
$ws.Range("E41").Value = "  -4.24%  "
$ws.Range("E42").Value = "  +0.32%  "
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("E44").Value = "  -4.69%  "
$ws.Range("E45").Value = "  -3.37%  "
$ws.Range("E46").Value = "  -3.69%  "
$ws.Range("E47").Value = "  -2.64%  "
$ws.Range("E48").Value = "  +0.76%  "
$ws.Range("E49").Value = "  -5.80%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("E50").Value = "  -0.19%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("E51").Value = "  -4.13%  "

Set-TextValue $ws.Range("D5") "590.46"
Set-TextValue $ws.Range("D6") "164.83"
Set-TextValue $ws.Range("D12") "0.448"
Set-TextValue $ws.Range("D13") "0.0000260"
Set-TextValue $ws.Range("D14") "35.95"
Set-TextValue $ws.Range("D18") "18.23"
Set-TextValue $ws.Range("D21") "10.66"
Set-TextValue $ws.Range("D22") "466.14"
Set-TextValue $ws.Range("D24") "82.67"
Set-TextValue $ws.Range("D25") "0.0000133"
Set-TextValue $ws.Range("D31") "2.75"
Set-TextValue $ws.Range("D32") "7.31"
Set-TextValue $ws.Range("D33") "2.22"
Set-TextValue $ws.Range("D34") "29.46"
Set-TextValue $ws.Range("D35") "9.00"
Set-TextValue $ws.Range("D40") "0.988"
Set-TextValue $ws.Range("D44") "0.303"
Set-TextValue $ws.Range("D45") "8.50"
Set-TextValue $ws.Range("D47") "45.20"
Set-TextValue $ws.Range("D48") "143.22"
Set-TextValue $ws.Range("D49") "384.28"
Set-TextValue $ws.Range("D50") "25.25"
Set-TextValue $ws.Range("D51") "0.0345"
